# Calibrate CDS params, sigmoid function for survivourship, liquidity risk
# step increase function, vol, and yield curve stress.

$wb = $excel.ActiveWorkbook

$wsRates  = $wb.Worksheets.Item("interest_rate_term_structure")
$wsCredit = $wb.Worksheets.Item("credit_risk_term_structure")
$wsRateP  = $wb.Worksheets.Item("interest_rate_parameters")
$wsCredP  = $wb.Worksheets.Item("credit_risk_parameters")

# ---------------------------------------------------------------------------
# interest_rate_parameters: recalibrated theta (B2) for the short-rate model
# ---------------------------------------------------------------------------
$wsRateP.Range("B2").Value = 0.031990400000000002

# ---------------------------------------------------------------------------
# credit_risk_parameters: re-fit sigmoid survivorship parameters and append
# the stepwise liquidity-risk / vol / yield-curve stress scenarios below the
# original calibration row.
# ---------------------------------------------------------------------------
$wsCredP.Range("A2").Value = 0.35590523000000002
$wsCredP.Range("B2").Value = 0.00149927
$wsCredP.Range("C2").Value = 0.0023764799999999998
$wsCredP.Range("D2").Value = 0.00197182

$wsCredP.Range("A3").Value = 0.060976130000000003
$wsCredP.Range("B3").Value = 0.0028177200000000001
$wsCredP.Range("C3").Value = 0.37078147
$wsCredP.Range("D3").Value = 0

$wsCredP.Range("A4").Value = 0.097805409999999995
$wsCredP.Range("B4").Value = 0.00057516999999999996
$wsCredP.Range("C4").Value = 0.12454273
$wsCredP.Range("D4").Value = 0.0029165900000000002

$wsCredP.Range("A5").Value = 0.35420099999999999
$wsCredP.Range("B5").Value = 0.00121853
$wsCredP.Range("C5").Value = 0.00238186
$wsCredP.Range("D5").Value = 0.00181

$wsCredP.Range("A6").Value = 0.37079495000000001
$wsCredP.Range("B6").Value = 0.0038134200000000001
$wsCredP.Range("C6").Value = 0.00250118
$wsCredP.Range("D6").Value = 0.0021209900000000001

$wsCredP.Range("A7").Value = 0.37072224999999998
$wsCredP.Range("B7").Value = 0.00309399
$wsCredP.Range("C7").Value = 0.0025071500000000001
$wsCredP.Range("D7").Value = 0.0036718499999999999

$wsCredP.Range("A8").Value = 0.36948254000000003
$wsCredP.Range("B8").Value = 0.00038684
$wsCredP.Range("C8").Value = 0.0024764700000000001
$wsCredP.Range("D8").Value = 0.0098715000000000001

# Row 9 carries the blp_amount-style scientific-notation formatting used
# elsewhere in the workbook (same numFmt as interest_rate_parameters!D2).
$wsCredP.Range("A9:D9").NumberFormat = "0.00E+00"
$wsCredP.Range("A9").Value = 0.63004337399999999
$wsCredP.Range("B9").Value = 0.00022447731499999999
$wsCredP.Range("C9").Value = 0.0041979461500000004
$wsCredP.Range("D9").Value = 0.0151173503

$wsCredP.Range("A10").Value = 0.42458847999999999
$wsCredP.Range("B10").Value = 0.0040863699999999998
$wsCredP.Range("C10").Value = 0.0028134100000000001
$wsCredP.Range("D10").Value = 0.0011898099999999999

$wsCredP.Range("A11").Value = 0.38846364
$wsCredP.Range("B11").Value = 0.00110464
$wsCredP.Range("C11").Value = 0.00259129
$wsCredP.Range("D11").Value = 0.001

$wsCredP.Range("A12").Value = 0.35590523000000002
$wsCredP.Range("B12").Value = 0.00149927
$wsCredP.Range("C12").Value = 0.0023764799999999998
$wsCredP.Range("D12").Value = 0.00197182

# Page setup picked up for the credit_risk_parameters sheet as its data grid
# grew.
$wsCredP.PageSetup.PaperSize = 9
$wsCredP.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selections / active-window bookkeeping, applied last so the final saved
# selection + active-sheet state matches what was left on screen.
# ---------------------------------------------------------------------------
$wsRates.Range("C2").Select()
$wsCredit.Range("D10").Select()
$wsRateP.Range("D2").Select()
$wsCredP.Range("A2:D2").Select()
